$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 9968.14
$ws.Range("B8").Value = 9926.4500000000007
$ws.Range("C8").Value = 80.11
$ws.Range("D8").Value = 79.77
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = -0.42
$ws.Range("G8").Value = 42609.488333333335
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"
$ws.Range("H8").Value = $true
